$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Recode column C (rows 2-45) from numeric 1/0 to text labels ---
# 1 -> "Позитивная оценка", 0 -> "Негативная оценка"
for ($r = 2; $r -le 45; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value()
    if ($v -eq 1) {
        $cell.Value = "Позитивная оценка"
    } else {
        $cell.Value = "Негативная оценка"
    }
}

# --- 2. Add two new rows (46, 47) for a new question block ---
$question15 = "15. Есть ли у вас функциональный руководитель? (не является руководителем по структуре, но ставит вам задачи)"

$ws.Cells.Item(46, 1).Value = $question15
$ws.Cells.Item(46, 2).Value = "Да"
$ws.Cells.Item(46, 3).Value = 1

$ws.Cells.Item(47, 1).Value = $question15
$ws.Cells.Item(47, 2).Value = "Нет"
$ws.Cells.Item(47, 3).Value = 0

# --- 3. Update the sheet view: scroll position + selection ---
$ws.Range("E46").Select()

# --- 4. Page setup (paper size / orientation) ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
